# ============================================================
# Edit: add "version" + "description" columns to the front of
# the "Export as TSV" sheet, add a new "version list" sheet,
# and shift existing header comments / data validations along.
# ============================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Export as TSV")

# --- Step 1: remove the existing header comments. They are anchored
# to fixed cell refs and will NOT move automatically when we insert
# new columns, so drop them now and re-add them at the correct,
# shifted locations afterwards.
$ws1.Range("A1").Comment.Delete()
$ws1.Range("B1").Comment.Delete()
$ws1.Range("C1").Comment.Delete()
$ws1.Range("D1").Comment.Delete()
$ws1.Range("E1").Comment.Delete()
$ws1.Range("F1").Comment.Delete()
$ws1.Range("G1").Comment.Delete()
$ws1.Range("H1").Comment.Delete()
$ws1.Range("I1").Comment.Delete()
$ws1.Range("J1").Comment.Delete()
$ws1.Range("K1").Comment.Delete()
$ws1.Range("L1").Comment.Delete()
$ws1.Range("M1").Comment.Delete()
$ws1.Range("N1").Comment.Delete()
$ws1.Range("O1").Comment.Delete()
$ws1.Range("P1").Comment.Delete()
$ws1.Range("Q1").Comment.Delete()
$ws1.Range("R1").Comment.Delete()
$ws1.Range("S1").Comment.Delete()
$ws1.Range("T1").Comment.Delete()
$ws1.Range("U1").Comment.Delete()
$ws1.Range("V1").Comment.Delete()
$ws1.Range("W1").Comment.Delete()
$ws1.Range("X1").Comment.Delete()
$ws1.Range("Y1").Comment.Delete()
$ws1.Range("Z1").Comment.Delete()
$ws1.Range("AA1").Comment.Delete()

# --- Step 2: insert two new blank columns at the front (A:B).
# This shifts all existing header cells, row data and data
# validation sqrefs from columns A..AA to C..AC automatically.
$ws1.Columns("A:B").Insert()

# --- Step 3: set the headers for the two new columns (and match the
# bold/centered/wrap-text style used by the rest of the header row).
$ws1.Range("A1").Value = "version"
$ws1.Range("B1").Value = "description"
$ws1.Range("A1:B1").Font.Bold = $true
$ws1.Range("A1:B1").HorizontalAlignment = -4108
$ws1.Range("A1:B1").WrapText = $true

# --- Step 4: re-create the header comments at their (shifted) cells.
$ws1.Range("A1").AddComment("Version of the schema to use when validating this metadata.")
$ws1.Range("B1").AddComment("Free-text description of this assay.")
$ws1.Range("C1").AddComment("HuBMAP Display ID of the donor of the assayed tissue.")
$ws1.Range("D1").AddComment("HuBMAP Display ID of the assayed tissue.")
$ws1.Range("E1").AddComment("Start date and time of assay, typically a date-time stamped folder generated by the acquisition instrument. YYYY-MM-DD hh:mm, where YYYY is the year, MM is the month with leading 0s, and DD is the day with leading 0s, hh is the hour with leading zeros, mm are the minutes with leading zeros.")
$ws1.Range("F1").AddComment("DOI for protocols.io referring to the protocol for this assay.")
$ws1.Range("G1").AddComment("Name of the person responsible for executing the assay.")
$ws1.Range("H1").AddComment("Email address for the operator.")
$ws1.Range("I1").AddComment("Name of the principal investigator responsible for the data.")
$ws1.Range("J1").AddComment("Email address for the principal investigator.")
$ws1.Range("K1").AddComment("Each assay is placed into one of the following 3 general categories: generation of images of microscopic entities, identification & quantitation of molecules by mass spectrometry, and determination of nucleotide sequence.")
$ws1.Range("L1").AddComment("The specific type of assay being executed.")
$ws1.Range("M1").AddComment("Analytes are the target molecules being measured with the assay.")
$ws1.Range("N1").AddComment("Specifies whether or not a specific molecule(s) is/are targeted for detection/measurement by the assay. The CODEX analyte is protein.")
$ws1.Range("O1").AddComment("An acquisition instrument is the device that contains the signal detection hardware and signal processing software. Assays generate signals such as light of various intensities or color or signals representing the molecular mass.")
$ws1.Range("P1").AddComment("Manufacturers of an acquisition instrument may offer various versions (models) of that instrument with different features or sensitivities. Differences in features or sensitivities may be relevant to processing or interpretation of the data.")
$ws1.Range("Q1").AddComment("Number of antibodies")
$ws1.Range("R1").AddComment("Number of fluorescent channels imaged during each cycle.")
$ws1.Range("S1").AddComment("Number of cycles of 1. oligo application, 2. fluor application, 3. dye inactivation.")
$ws1.Range("T1").AddComment("the total number of acquisitions performed on microscope to collect autofluorescence/background or stained signal.")
$ws1.Range("U1").AddComment("The width of a pixel.")
$ws1.Range("V1").AddComment("The unit of measurement of the width of a pixel.")
$ws1.Range("W1").AddComment("The height of a pixel")
$ws1.Range("X1").AddComment("The unit of measurement of the height of a pixel.")
$ws1.Range("Y1").AddComment("DOI for analysis protocols.io for this assay.")
$ws1.Range("Z1").AddComment("DOI for protocols.io for the overall process.")
$ws1.Range("AA1").AddComment("Relative path to file with antibody information for this dataset.")
$ws1.Range("AB1").AddComment("Relative path to file with ORCID IDs for contributors for this dataset.")
$ws1.Range("AC1").AddComment("Relative path to file or directory with instrument data. Downstream processing will depend on filename extension conventions.")

# --- Step 5: add data validation (dropdown) for the new "version" column.
$rngVersion = $ws1.Range("A2:A1048576")
$rngVersion.Validation.Add(3, 1, 1, "='version list'!`$A`$1:`$A`$1")
$rngVersion.Validation.ErrorTitle = "Value must come from list"
$rngVersion.Validation.ErrorMessage = "Value must be one of: 1."
$rngVersion.Validation.ShowInput = $true
$rngVersion.Validation.ShowError = $true

# --- Step 6: add the new "version list" sheet right after "Export as TSV".
$versionSheet = $wb.Worksheets.Add($null, $ws1)
$versionSheet.Name = "version list"
$versionSheet.Range("A1").Value = "1"

Write-Output "edit complete"
